# Fixing and update digisales new version
# Adds a new scenario row (SCD0206-SCD0209) to the "Scenario" sheet,
# mirroring the formatting of the row above it, and moves the sheet's
# view/selection down to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario")

# New scenario / note content (row 11)
$ws.Range("A11").Value = "SCD0206-SCD0209"
$ws.Range("B11").Value = "1. Button Export Excel nya tidak tersedia`n2. Step Export Excel tidak dapat dilakukan"

# Match formatting used by the other data rows:
#  - column A: vertically centered (style index 5 on existing rows)
#  - column B: wrapped text (style index 1 on existing rows)
$ws.Range("A11").VerticalAlignment = -4108  # xlCenter
$ws.Range("B11").WrapText = $true

# Row height matches the shorter, 2-line note text
$ws.Rows.Item(11).RowHeight = 30

# Move the view so the new row is visible and select it
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B11").Select()
